$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Commit: "Cambio di segno alle derivate di controllo del latero direzionale"
# The sign of the control derivatives (w.r.t. delta_r in column F and
# delta_a in column G) for the lateral-directional coefficients
# CY, CL, CN (rows 10-15) is flipped.
$ws.Range("F10").Value = 70.884857142857157
$ws.Range("G10").Value = -61.474867924528297

$ws.Range("F11").Value = 17.120822454308097
$ws.Range("G11").Value = -109.97231746031748

$ws.Range("F12").Value = -250.18367139959437
$ws.Range("G12").Value = -16.1869387755102

$ws.Range("F13").Value = 35.937440944881885
$ws.Range("G13").Value = -2.9880314960629928

$ws.Range("F14").Value = 83.03658333333334
$ws.Range("G14").Value = -57.816196581196586

$ws.Range("F15").Value = 20.683461538461533
$ws.Range("G15").Value = -107.42181203007519
